$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.664.73"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "'1.896.37"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'310.21"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").Value = "'0.9990"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +6.12%  "
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").Value = "'0.07237"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").Value = "'21.08"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "'0.9010"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").Value = "'1.891.58"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "'0.07626"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "'5.432"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "'91.66"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "'0.9999"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "'0.000008662"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "'14.31"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "'0.9992"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'27.696.89"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "'2.127.13"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").Value = "'10.80"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'6.602"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").Value = "'153.06"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("D28").Value = "'2.172"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").Value = "'113.96"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").Value = "'4.822"
$ws.Range("E30").Value = "  -2.21%  "
$ws.Range("D31").Value = "'4.817"
$ws.Range("E31").Value = "  +3.73%  "
$ws.Range("D32").Value = "'0.09156"
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("D33").Value = "'0.05265"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").Value = "'3.156"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("D35").Value = "'1.222"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").Value = "'0.7720"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "'0.02082"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").Value = "'2.560"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("E39").Value = "  +2.22%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5559"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.089"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").Value = "'6.685"
$ws.Range("E42").Value = "  -4.44%  "
$ws.Range("D43").Value = "'117.46"
$ws.Range("E43").Value = "  +5.55%  "
$ws.Range("D44").Value = "'8.705"
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").Value = "'0.4795"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("D48").Value = "'0.9984"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").Value = "'1.590"
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("D51").Value = "'37.00"
$ws.Range("E51").Value = "  -0.13%  "
